$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.187.98"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "1.641.58"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "217.18"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "0.526"
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Value = "19.97"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "1.872.68"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "1.660.15"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").Value = "67.13"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "27.190.17"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").Value = "218.73"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "6.95"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("E23").Value = "  +3.44%  "
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").Value = "147.54"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "7.44"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").Value = "1.303.74"
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("D37").Value = "0.0176"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "0.550"
$ws.Range("E38").Value = "  +3.24%  "
$ws.Range("E39").Value = "  +3.09%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  +5.98%  "
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("D44").Value = "1.782.20"
$ws.Range("D45").Value = "61.76"
$ws.Range("D46").Value = "91.74"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "7.66"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("E51").Value = "  +0.53%  "
